# Druid archetype sheet - "First two druid subclasses completed"
#
# 1. Novice tier (row 7, Level 6, Discipline choice) gains its descriptive
#    blurb in column H, explaining the 3-of-6 discipline choice.
# 2. The Expert tier's "List2 Name" entry (row 14, D column) is renamed from
#    the placeholder "Ancient Knowledge" to "Ancient Powers".
# 3. The stray/duplicate "Aspect Feature VI" bonus entry on row 19 is removed
#    (Master tier no longer lists it here).
# 4. Selection cursor ends up parked on A20 (where the author's cursor was
#    left after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New discipline-choice description for the Novice tier, column H (row 7).
$ws.Range("H7").Value = "Choose 3 from Elemental, Telepathy, Bewitchment, Healing, Warding and Alteration. "
$ws.Range("H7").WrapText = $true
$ws.Range("H7").HorizontalAlignment = -4152
$ws.Rows(7).RowHeight = 46.25

# 2. Rename "Ancient Knowledge" -> "Ancient Powers" (Expert tier, List2 Name).
$ws.Range("D14").Value = "Ancient Powers"

# 3. Remove the extraneous "Aspect Feature VI" bonus text on row 19.
$ws.Range("D19").ClearContents()

# 4. Leave the selection where the author left it.
$ws.Range("A20").Select() | Out-Null
